$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.205.12'
$ws.Range('D2').Style = $cellStyle
$ws.Range('E2').Value = '  +1.94%  '
$cellStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.093.87'
$ws.Range('D3').Style = $cellStyle
$ws.Range('E3').Value = '  +3.49%  '
$ws.Range('E4').Value = '  +0.17%  '
$cellStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.59'
$ws.Range('D5').Style = $cellStyle
$ws.Range('E5').Value = '  +1.44%  '
$cellStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.659'
$ws.Range('D6').Style = $cellStyle
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  +0.09%  '
$cellStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '52.12'
$ws.Range('D8').Style = $cellStyle
$ws.Range('E8').Value = '  +14.79%  '
$cellStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '61.57'
$ws.Range('D9').Style = $cellStyle
$ws.Range('E9').Value = '  +3.05%  '
$cellStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.373'
$ws.Range('D10').Style = $cellStyle
$ws.Range('E10').Value = '  +2.10%  '
$ws.Range('E11').Value = '  +3.33%  '
$ws.Range('E12').Value = '  +7.35%  '
$cellStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.10'
$ws.Range('D13').Style = $cellStyle
$ws.Range('E13').Value = '  +2.51%  '
$cellStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.397.38'
$ws.Range('D14').Style = $cellStyle
$ws.Range('E14').Value = '  +3.54%  '
$ws.Range('E15').Value = '  +2.41%  '
$cellStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.091.86'
$ws.Range('D16').Style = $cellStyle
$ws.Range('E16').Value = '  +3.50%  '
$cellStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.12'
$ws.Range('D17').Style = $cellStyle
$ws.Range('E17').Value = '  +3.72%  '
$cellStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.175.56'
$ws.Range('D18').Style = $cellStyle
$ws.Range('E18').Value = '  +1.70%  '
$cellStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.27'
$ws.Range('D19').Style = $cellStyle
$ws.Range('E19').Value = '  +1.31%  '
$cellStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.11'
$ws.Range('D20').Style = $cellStyle
$ws.Range('E20').Value = '  +8.50%  '
$cellStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0839'
$ws.Range('D21').Style = $cellStyle
$ws.Range('E21').Value = '  +1.98%  '
$cellStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '240.36'
$ws.Range('D22').Style = $cellStyle
$ws.Range('E22').Value = '  +1.51%  '
$cellStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.22'
$ws.Range('D23').Style = $cellStyle
$ws.Range('E23').Value = '  +5.98%  '
$ws.Range('E24').Value = '  +0.00%  '
$cellStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.46'
$ws.Range('D25').Style = $cellStyle
$cellStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.38'
$ws.Range('D26').Style = $cellStyle
$ws.Range('E26').Value = '  +4.59%  '
$ws.Range('E27').Value = '  +6.61%  '
$cellStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.70'
$ws.Range('D28').Style = $cellStyle
$ws.Range('E28').Value = '  +3.37%  '
$ws.Range('E29').Value = '  +1.70%  '
$cellStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.123'
$ws.Range('D30').Style = $cellStyle
$ws.Range('E30').Value = '  +0.44%  '
$cellStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.74'
$ws.Range('D31').Style = $cellStyle
$ws.Range('E31').Value = '  +4.95%  '
$ws.Range('E32').Value = '  +25.46%  '
$cellStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.49'
$ws.Range('D33').Style = $cellStyle
$ws.Range('E33').Value = '  +1.54%  '
$cellStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0610'
$ws.Range('D34').Style = $cellStyle
$ws.Range('E34').Value = '  +2.62%  '
$cellStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0909'
$ws.Range('D35').Style = $cellStyle
$ws.Range('E35').Value = '  +10.67%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cellStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.87'
$ws.Range('D37').Style = $cellStyle
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cellStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.29'
$ws.Range('D38').Style = $cellStyle
$ws.Range('E38').Value = '  +7.37%  '
$ws.Range('E39').Value = '  +1.27%  '
$ws.Range('E40').Value = '  -0.78%  '
$cellStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.26'
$ws.Range('D41').Style = $cellStyle
$ws.Range('E41').Value = '  +11.81%  '
$ws.Range('E42').Value = '  +3.34%  '
$ws.Range('E43').Value = '  +4.70%  '
$ws.Range('E44').Value = '  +2.08%  '
$cellStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0919'
$ws.Range('D45').Style = $cellStyle
$ws.Range('E45').Value = '  +13.34%  '
$cellStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.74'
$ws.Range('D46').Style = $cellStyle
$ws.Range('E46').Value = '  -0.31%  '
$cellStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.319.48'
$ws.Range('D47').Style = $cellStyle
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('E48').Value = '  +6.68%  '
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$cellStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.83'
$ws.Range('D49').Style = $cellStyle
$ws.Range('E49').Value = '  +83.30%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cellStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.97'
$ws.Range('D50').Style = $cellStyle
$ws.Range('E50').Value = '  +13.15%  '
$cellStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.286.48'
$ws.Range('D51').Style = $cellStyle
$ws.Range('E51').Value = '  +2.70%  '
